$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 88

# Numeric / reused-string cells first (order irrelevant for the shared
# string table since these values already exist in it).
$ws.Cells.Item($row, 1).Value = 43224.375
$ws.Cells.Item($row, 2).Value = "分类O"
$ws.Cells.Item($row, 3).Value = "14分类"
$ws.Cells.Item($row, 4).Value = "batch_size=100 low_nums=2 use_biases=yes use_bn_input=True   dropout_low=0.8 "
$ws.Cells.Item($row, 6).Value = 0.65
$ws.Cells.Item($row, 7).Value = 0.65
$ws.Cells.Item($row, 8).Value = 0.98
$ws.Cells.Item($row, 9).Value = 0.98

# New unique strings must be entered in the same order the author typed
# them so the shared-string table indices line up with the target file:
# J (result) -> E (data) -> K (command) -> L (logs).
$ws.Cells.Item($row, 10).Value = "经过约23小时，较长时间泛化精度都在拟合精度之上，后面拟合开始超越，泛化波动较大，最好达到0.66。"
$ws.Cells.Item($row, 5).Value = "最高标签，重新训练，经过归一化数据加PCA及Wavelet处理数据train-hjxh365-2018-4-16-day-high-norm-pca99-wavelet20"
$ws.Cells.Item($row, 11).Value = "python feed_run.py --output_mode=classes --output_nodes=14 --input_nums=96 --input_nodes=96 --low_nums=2 --low_nodes=96 --low_fun=elu --one_hot=True --input_fun=elu --batch_size=100 --learning_rate=0.001 --train_mode=Adadelta --eval_size=5400 --test_size=1339 --use_biases=yes  --use_bn_input=True --dropout_low=0.8"
$ws.Cells.Item($row, 12).Value = "logs-hjxh-2018-5-4-high-norm-pca99-wavelet20-percent65"

$ws.Rows.Item($row).RowHeight = 82.5

$ws.Range("E86").Select()
